# website update 12-5-2017 &logboek
# Adds three new logboek entries (rows 35, 37 and 38) to the worksheet,
# mirroring the style/format of the existing rows above them.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a numeric-looking value as TEXT (so it is stored as a
# shared string rather than being auto-converted into a number), without
# leaving behind any extra number-format/style definitions. We do this by
# entering it as a quoted-string formula and then collapsing the formula
# down to its static result via PasteSpecial(xlPasteValues).
function Set-TextValue($cell, [string]$text) {
    $cell.Formula = '="' + $text + '"'
    $cell.Copy() | Out-Null
    $cell.PasteSpecial(-4163) | Out-Null  # xlPasteValues
}

# Helper: write a date value into a cell while re-using the same
# number-format/style as the other date cells in column A (style index 3
# in the original workbook, built from copying cell A34's formatting).
function Set-DateValue($cell, $serial) {
    $ws.Cells.Item(34, 1).Copy() | Out-Null
    $cell.PasteSpecial(-4122) | Out-Null  # xlPasteFormats
    $cell.Value = $serial
}

# ---- Row 35 : 11-5-2017 --------------------------------------------------
Set-DateValue $ws.Cells.Item(35, 1) 42866
Set-TextValue $ws.Cells.Item(35, 3) "23.30"
$ws.Cells.Item(35, 5).Value = "thuis"
$ws.Cells.Item(35, 7).Value = "werken aan de php van resultaten.php."

# ---- Row 37 : 12-5-2017 ---------------------------------------------------
Set-DateValue $ws.Cells.Item(37, 1) 42867
Set-TextValue $ws.Cells.Item(37, 3) "8.40"
$ws.Cells.Item(37, 5).Value = "School"
$ws.Cells.Item(37, 7).Value = "werken aan de html van invoeren teams en spelers"

# ---- Row 38 : 12-5-2017 ---------------------------------------------------
Set-DateValue $ws.Cells.Item(38, 1) 42867
Set-TextValue $ws.Cells.Item(38, 3) "9.30"
$ws.Cells.Item(38, 5).Value = "School"
$ws.Cells.Item(38, 7).Value = "werken aan de php van invoeren teams en spelers"

# ---- Update the view so the new rows are visible/selected ----------------
$ws.Activate() | Out-Null
$ws.Range("G38").Select() | Out-Null
$excel.ActiveWindow.ScrollRow = 34
$excel.ActiveWindow.ScrollColumn = 1
